$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# New scenario: "repaymentstrategy" value changes from "Mifos style" to the
# new "Penalties, Fees, Interest, Principal order" strategy, with the cell
# re-styled to left/top aligned text.
$cell = $wsInput.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

# The ProductLoanInput sheet becomes the active/selected sheet & cell,
# while ProductLoanOutput loses the tab-selected state it previously had.
$wsInput.Activate() | Out-Null
$wsInput.Range("B17").Select() | Out-Null
